$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B8").Value = "x"
$ws.Range("B13").Value = "x"
$ws.Range("B14").Value = "x"

$ws.Range("B3").Select()
